$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Treatment Technologies"
$ws.Cells.Item(2, 2).Value = "JWC Environmental"
$ws.Cells.Item(2, 3).Value = "Muffin Monster"
$ws.Cells.Item(2, 4).Value = "Screening"

$ws.Cells.Item(3, 1).Value = "Treatment Technologies"
$ws.Cells.Item(3, 2).Value = "JWC Environmental"
$ws.Cells.Item(3, 3).Value = "Channel Monster"
$ws.Cells.Item(3, 4).Value = "Screening"

$ws.Cells.Item(4, 1).Value = "Treatment Technologies"
$ws.Cells.Item(4, 2).Value = "JWC Environmental"
$ws.Cells.Item(4, 3).Value = "Auger Monster"
$ws.Cells.Item(4, 4).Value = "Screening"

$ws.Cells.Item(5, 1).Value = "Treatment Technologies"
$ws.Cells.Item(5, 2).Value = "JWC Environmental"
$ws.Cells.Item(5, 3).Value = "Fine Screens"
$ws.Cells.Item(5, 4).Value = "Screening"

$ws.Cells.Item(6, 1).Value = "Treatment Technologies"
$ws.Cells.Item(6, 2).Value = "JWC Environmental"
$ws.Cells.Item(6, 3).Value = "Band Screens"
$ws.Cells.Item(6, 4).Value = "Screening"

$ws.Cells.Item(7, 1).Value = "Treatment Technologies"
$ws.Cells.Item(7, 2).Value = "JWC Environmental"
$ws.Cells.Item(7, 3).Value = "Drum Screens"
$ws.Cells.Item(7, 4).Value = "Screening"

$ws.Cells.Item(8, 1).Value = "Treatment Technologies"
$ws.Cells.Item(8, 2).Value = "JWC Environmental"
$ws.Cells.Item(8, 3).Value = "Strom Screens"
$ws.Cells.Item(8, 4).Value = "Screening"

$ws.Cells.Item(9, 1).Value = "Treatment Technologies"
$ws.Cells.Item(9, 2).Value = "JWC Environmental"
$ws.Cells.Item(9, 3).Value = "Screening Washer Monster"
$ws.Cells.Item(9, 4).Value = "Screening"

$ws.Cells.Item(10, 1).Value = "Treatment Technologies"
$ws.Cells.Item(10, 2).Value = "JWC Environmental"
$ws.Cells.Item(10, 3).Value = "Honey Monster Septage Receiving Station"
$ws.Cells.Item(10, 4).Value = "Grease and FOG Removal"

$ws.Cells.Item(11, 1).Value = "Treatment Technologies"
$ws.Cells.Item(11, 2).Value = "IPEC"
$ws.Cells.Item(11, 3).Value = "Rotary Drum Sludge Thickeners"
$ws.Cells.Item(11, 4).Value = "Clarification"

$ws.Cells.Item(12, 1).Value = "Treatment Technologies"
$ws.Cells.Item(12, 2).Value = "IPEC"
$ws.Cells.Item(12, 3).Value = "Internally-Fed Rotary Screens"
$ws.Cells.Item(12, 4).Value = "Screening"

$ws.Cells.Item(13, 1).Value = "Treatment Technologies"
$ws.Cells.Item(13, 2).Value = "IPEC"
$ws.Cells.Item(13, 3).Value = "Static Screens"
$ws.Cells.Item(13, 4).Value = "Screening"

$ws.Cells.Item(14, 1).Value = "Treatment Technologies"
$ws.Cells.Item(14, 2).Value = "FRC Systems International"
$ws.Cells.Item(14, 3).Value = "Dissolved Air Floatation (DAF)"
$ws.Cells.Item(14, 4).Value = "Clarification"

$ws.Cells.Item(15, 1).Value = "Treatment Technologies"
$ws.Cells.Item(15, 2).Value = "V-Fold"
$ws.Cells.Item(15, 3).Value = "Sludge Dewatering"
$ws.Cells.Item(15, 4).Value = "Filtration"

$ws.Cells.Item(16, 1).Value = "Treatment Technologies"
$ws.Cells.Item(16, 2).Value = "Chemco Systems"
$ws.Cells.Item(16, 3).Value = "Powder Activated Carbon"
$ws.Cells.Item(16, 4).Value = "Chemical Feed"

$ws.Cells.Item(17, 1).Value = "Treatment Technologies"
$ws.Cells.Item(17, 2).Value = "Chemco Systems"
$ws.Cells.Item(17, 3).Value = "Lime Slacker"
$ws.Cells.Item(17, 4).Value = "Chemical Feed"

$ws.Cells.Item(18, 1).Value = "Treatment Technologies"
$ws.Cells.Item(18, 2).Value = "Chemco Systems"
$ws.Cells.Item(18, 3).Value = "Hydrated Lime"
$ws.Cells.Item(18, 4).Value = "Chemical Feed"

$ws.Cells.Item(19, 1).Value = "Treatment Technologies"
$ws.Cells.Item(19, 2).Value = "Chemco Systems"
$ws.Cells.Item(19, 3).Value = "Soda Ash"
$ws.Cells.Item(19, 4).Value = "Chemical Feed"

$ws.Cells.Item(20, 1).Value = "Treatment Technologies"
$ws.Cells.Item(20, 2).Value = "Chemco Systems"
$ws.Cells.Item(20, 3).Value = "Bulk Bag Unloader"
$ws.Cells.Item(20, 4).Value = "Chemical Feed"

$ws.Cells.Item(21, 1).Value = "Treatment Technologies"
$ws.Cells.Item(21, 2).Value = "Chemco Systems"
$ws.Cells.Item(21, 3).Value = "Silo Systems"
$ws.Cells.Item(21, 4).Value = "Chemical Feed"

$ws.Cells.Item(22, 1).Value = "Treatment Technologies"
$ws.Cells.Item(22, 2).Value = "Clearas"
$ws.Cells.Item(22, 3).Value = "Advanced Biological Nutrient Recovery (ABNR) technology"
$ws.Cells.Item(22, 4).Value = "Nutrient Recovery"

$ws.Cells.Item(23, 1).Value = "Treatment Technologies"
$ws.Cells.Item(23, 2).Value = "Fibracast"
$ws.Cells.Item(23, 3).Value = "FibrePlate Membrane Bioreactor (MBR) Technology"
$ws.Cells.Item(23, 4).Value = "Filtration"

$ws.Cells.Item(24, 1).Value = "Treatment Technologies"
$ws.Cells.Item(24, 2).Value = "Kruger"
$ws.Cells.Item(24, 3).Value = "Phased Oxidation Ditches"
$ws.Cells.Item(24, 4).Value = "Aeration"

$ws.Cells.Item(25, 1).Value = "Treatment Technologies"
$ws.Cells.Item(25, 2).Value = "Kruger"
$ws.Cells.Item(25, 3).Value = "Discfilter"
$ws.Cells.Item(25, 4).Value = "Filtration"

$ws.Cells.Item(26, 1).Value = "Treatment Technologies"
$ws.Cells.Item(26, 2).Value = "Kruger"
$ws.Cells.Item(26, 3).Value = "BioCon Thermal Dryer"
$ws.Cells.Item(26, 4).Value = "Sludge Management"

$ws.Cells.Item(27, 1).Value = "Treatment Technologies"
$ws.Cells.Item(27, 2).Value = "Kruger"
$ws.Cells.Item(27, 3).Value = "ACTIFLO"
$ws.Cells.Item(27, 4).Value = "Clarification"

$ws.Cells.Item(28, 1).Value = "Treatment Technologies"
$ws.Cells.Item(28, 2).Value = "DTE Environmental"
$ws.Cells.Item(28, 3).Value = "Grit Classifier"
$ws.Cells.Item(28, 4).Value = "Screening"

$ws.Cells.Item(29, 1).Value = "Treatment Technologies"
$ws.Cells.Item(29, 2).Value = "Rodney Hunt"
$ws.Cells.Item(29, 3).Value = "Sluice Gates"
$ws.Cells.Item(29, 4).Value = "Flow Control"

$ws.Cells.Item(30, 1).Value = "Treatment Technologies"
$ws.Cells.Item(30, 2).Value = "Rodney Hunt"
$ws.Cells.Item(30, 3).Value = "Slide Gates"
$ws.Cells.Item(30, 4).Value = "Flow Control"

$ws.Cells.Item(31, 1).Value = "Treatment Technologies"
$ws.Cells.Item(31, 2).Value = "Entek Technologies"
$ws.Cells.Item(31, 3).Value = "Aeration Systems"
$ws.Cells.Item(31, 4).Value = "Aeration"

$ws.Cells.Item(32, 1).Value = "Treatment Technologies"
$ws.Cells.Item(32, 2).Value = "BNR Systems"
$ws.Cells.Item(32, 3).Value = "Shaftless Spiral Conveyors"
$ws.Cells.Item(32, 4).Value = "Sludge Management"

$ws.Cells.Item(33, 1).Value = "Treatment Technologies"
$ws.Cells.Item(33, 2).Value = "BNR Systems"
$ws.Cells.Item(33, 3).Value = "Live Bottom Hoppers"
$ws.Cells.Item(33, 4).Value = "Sludge Management"

$ws.Cells.Item(34, 1).Value = "Treatment Technologies"
$ws.Cells.Item(34, 2).Value = "BNR Systems"
$ws.Cells.Item(34, 3).Value = "Chain and Rake Screens"
$ws.Cells.Item(34, 4).Value = "Screening"

$ws.Cells.Item(35, 1).Value = "Treatment Technologies"
$ws.Cells.Item(35, 2).Value = "BNR Systems"
$ws.Cells.Item(35, 3).Value = "Fine Screens"
$ws.Cells.Item(35, 4).Value = "Screening"

$ws.Cells.Item(36, 1).Value = "Treatment Technologies"
$ws.Cells.Item(36, 2).Value = "BNR Systems"
$ws.Cells.Item(36, 3).Value = "Packages Headwork Systems for Screening and Grit Removal"
$ws.Cells.Item(36, 4).Value = "Screening"

$ws.Cells.Item(37, 1).Value = "Treatment Technologies"
$ws.Cells.Item(37, 2).Value = "BNR Systems"
$ws.Cells.Item(37, 3).Value = "Grit Vortek"
$ws.Cells.Item(37, 4).Value = "Screening"

$ws.Cells.Item(38, 1).Value = "Treatment Technologies"
$ws.Cells.Item(38, 2).Value = "BNR Systems"
$ws.Cells.Item(38, 3).Value = "Grit Classifier"
$ws.Cells.Item(38, 4).Value = "Screening"

$ws.Cells.Item(39, 1).Value = "Treatment Technologies"
$ws.Cells.Item(39, 2).Value = "BNR Systems"
$ws.Cells.Item(39, 3).Value = "Screenings Washer"
$ws.Cells.Item(39, 4).Value = "Screening"

